$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 166683860
$ws.Range("I80").Value = 200000620
$ws.Range("J80").Value = 100003
$ws.Range("K80").Value = 600001860
$ws.Range("L80").Value = 300009
$ws.Range("M80").Value = -600000862
$ws.Range("N80").Value = -302005
$ws.Range("H83").Value = 166683860
$ws.Range("I83").Value = 200000620
$ws.Range("J83").Value = 100003
$ws.Range("K83").Value = 1800005580
$ws.Range("L83").Value = 900027
$ws.Range("M83").Value = -1800000588
$ws.Range("N83").Value = -910011
$ws.Range("H106").Value = 4169172.2
$ws.Range("I106").Value = 4447003.5
$ws.Range("K106").Value = 4447003.5
$ws.Range("M106").Value = -4446372.5
$ws.Range("H125").Value = 3982.8572
$ws.Range("J125").Value = 5286
$ws.Range("L125").Value = 47574
$ws.Range("N125").Value = -52494
$ws.Range("H135").Value = 1847.4348
$ws.Range("I135").Value = 1449.8235
$ws.Range("J135").Value = 2974
$ws.Range("K135").Value = 13048.4115
$ws.Range("L135").Value = 26766
$ws.Range("M135").Value = -10513.4115
$ws.Range("N135").Value = -31836
$ws.Range("H137").Value = 5557426
$ws.Range("I137").Value = 1397.5
$ws.Range("K137").Value = 4192.5
$ws.Range("M137").Value = -1642.5
$ws.Range("H141").Value = 5455.2354
$ws.Range("I141").Value = 5455.2354
$ws.Range("K141").Value = 16365.7062
$ws.Range("M141").Value = -11185.7062
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 262450
$ws.Range("I34").Value = 24900
$ws.Range("K34").Value = 24900
$ws.Range("M34").Value = -24629
$ws.Range("H44").Value = 69000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H55").Value = 10048
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 8686.423000000001
$ws.Range("I61").Value = 17322.223
$ws.Range("K61").Value = 17322.223
$ws.Range("M61").Value = -17110.223
$ws.Range("H63").Value = 5550
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 5550
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("M63").Value = 5550
$ws.Range("N63").Value = -6922
$ws.Range("H66").Value = 5550
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 5550
$ws.Range("K66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("M66").Value = 27750
$ws.Range("N66").Value = -34614
$ws.Range("H132").Value = 4166.185
$ws.Range("I132").Value = 4561.4375
$ws.Range("J132").Value = 3591.2727
$ws.Range("K132").Value = 13684.3125
$ws.Range("L132").Value = 10773.8181
$ws.Range("M132").Value = -11154.3125
$ws.Range("N132").Value = -15833.8181
$ws.Range("H136").Value = 8686.423000000001
$ws.Range("I136").Value = 17322.223
$ws.Range("K136").Value = 51966.66900000001
$ws.Range("M136").Value = -49416.66900000001
$ws.Range("H141").Value = 167704.25
$ws.Range("J141").Value = 156942.33
$ws.Range("L141").Value = 156942.33
$ws.Range("N141").Value = -167302.33
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2338.75
$ws.Range("J99").Value = 3525
$ws.Range("L99").Value = 3525
$ws.Range("N99").Value = -6521
$ws.Range("H105").Value = 4966.2
$ws.Range("I105").Value = 6666.5
$ws.Range("K105").Value = 6666.5
$ws.Range("M105").Value = -4919.5
$ws.Range("H133").Value = 97472.75
$ws.Range("J133").Value = 97472.75
$ws.Range("L133").Value = 97472.75
$ws.Range("N133").Value = -107592.75
$ws.Range("H134").Value = 4875.353
$ws.Range("I134").Value = 4736.375
$ws.Range("J134").Value = 4998.8887
$ws.Range("K134").Value = 14209.125
$ws.Range("L134").Value = 14996.6661
$ws.Range("M134").Value = -11674.125
$ws.Range("N134").Value = -20066.6661
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1789.6666
$ws.Range("I16").Value = 1789.6666
$ws.Range("K16").Value = 1789.6666
$ws.Range("M16").Value = -1502.6666
$ws.Range("H31").Value = 6343.409
$ws.Range("I31").Value = 2397.8572
$ws.Range("J31").Value = 8184.6665
$ws.Range("K31").Value = 2397.8572
$ws.Range("L31").Value = 8184.6665
$ws.Range("M31").Value = -2102.8572
$ws.Range("N31").Value = -8774.666499999999
$ws.Range("H34").Value = 6343.409
$ws.Range("I34").Value = 2397.8572
$ws.Range("J34").Value = 8184.6665
$ws.Range("K34").Value = 2397.8572
$ws.Range("L34").Value = 8184.6665
$ws.Range("M34").Value = -2195.8572
$ws.Range("N34").Value = -8588.666499999999
$ws.Range("H113").Value = 1789.6666
$ws.Range("I113").Value = 1789.6666
$ws.Range("K113").Value = 1789.6666
$ws.Range("M113").Value = 380.3334
$ws.Range("H134").Value = 1751.5
$ws.Range("I134").Value = 1751.5
$ws.Range("K134").Value = 5254.5
$ws.Range("M134").Value = -2719.5
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("N140").Value = 0
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8538.166999999999
$ws.Range("I3").Value = 8538.166999999999
$ws.Range("K3").Value = 25614.501
$ws.Range("M3").Value = -25502.501
$ws.Range("H39").Value = 4040.5
$ws.Range("J39").Value = 3837.25
$ws.Range("L39").Value = 11511.75
$ws.Range("N39").Value = -12099.75
$ws.Range("H122").Value = 524.3125
$ws.Range("I122").Value = 414.66666
$ws.Range("J122").Value = 590.1
$ws.Range("K122").Value = 3731.99994
$ws.Range("L122").Value = 5310.900000000001
$ws.Range("M122").Value = -1281.99994
$ws.Range("N122").Value = -10210.9
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2612.7
$ws.Range("I126").Value = 2368
$ws.Range("J126").Value = 2744.4614
$ws.Range("K126").Value = 7104
$ws.Range("L126").Value = 8233.3842
$ws.Range("M126").Value = -4634
$ws.Range("N126").Value = -13173.3842
$ws.Range("H132").Value = 3999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").Value = 11997
$ws.Range("N132").Value = -17057
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2177.15
$ws.Range("I22").Value = 1951.5
$ws.Range("J22").Value = 2515.625
$ws.Range("K22").Value = 1951.5
$ws.Range("L22").Value = 2515.625
$ws.Range("M22").Value = -1656.5
$ws.Range("N22").Value = -3105.625
$ws.Range("H27").Value = 2177.15
$ws.Range("I27").Value = 1951.5
$ws.Range("J27").Value = 2515.625
$ws.Range("K27").Value = 1951.5
$ws.Range("L27").Value = 2515.625
$ws.Range("M27").Value = -1844.5
$ws.Range("N27").Value = -2729.625
$ws.Range("H55").Value = 1902.2778
$ws.Range("I55").Value = 2138.3333
$ws.Range("K55").Value = 2138.3333
$ws.Range("M55").Value = -1965.3333
$ws.Range("H61").Value = 3663.4092
$ws.Range("I61").Value = 1278.3334
$ws.Range("J61").Value = 8774.286
$ws.Range("K61").Value = 1278.3334
$ws.Range("L61").Value = 8774.286
$ws.Range("M61").Value = -1076.3334
$ws.Range("N61").Value = -9178.286
$ws.Range("H113").Value = 3663.4092
$ws.Range("I113").Value = 1278.3334
$ws.Range("J113").Value = 8774.286
$ws.Range("K113").Value = 1278.3334
$ws.Range("L113").Value = 8774.286
$ws.Range("M113").Value = 891.6666
$ws.Range("N113").Value = -13114.286
$ws.Range("H132").Value = 3046.0908
$ws.Range("I132").Value = 2401.3333
$ws.Range("J132").Value = 3287.875
$ws.Range("K132").Value = 7203.999899999999
$ws.Range("L132").Value = 9863.625
$ws.Range("M132").Value = -4673.999899999999
$ws.Range("N132").Value = -14923.625
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 39188.8
$ws.Range("J74").Value = 39188.8
$ws.Range("L74").Value = 39188.8
$ws.Range("N74").Value = -41060.8
$ws.Range("H77").Value = 39188.8
$ws.Range("J77").Value = 39188.8
$ws.Range("L77").Value = 117566.4
$ws.Range("N77").Value = -126926.4
$ws.Range("H81").Value = 3161.6667
$ws.Range("I81").Value = 1992.1
$ws.Range("J81").Value = 4623.625
$ws.Range("K81").Value = 3984.2
$ws.Range("L81").Value = 9247.25
$ws.Range("M81").Value = -2923.2
$ws.Range("N81").Value = -11369.25
$ws.Range("H84").Value = 3161.6667
$ws.Range("I84").Value = 1992.1
$ws.Range("J84").Value = 4623.625
$ws.Range("K84").Value = 19921
$ws.Range("L84").Value = 46236.25
$ws.Range("M84").Value = -14617
$ws.Range("N84").Value = -56844.25
$ws.Range("H132").Value = 3211.5833
$ws.Range("I132").Value = 2995.7144
$ws.Range("J132").Value = 3513.8
$ws.Range("K132").Value = 8987.143199999999
$ws.Range("L132").Value = 10541.4
$ws.Range("M132").Value = -6457.143199999999
$ws.Range("N132").Value = -15601.4
